$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Split "EmailAddress" into "Ema" / "ilAddress" with a _GoBack bookmark
#    placed between the two runs (simulating the cursor position after an
#    edit at that spot).
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$splitPos = $p6.Range.Start + 3
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Wrap "TimeStamp" with spellcheck proofErr markers inside the
#    "TransactionId, TimeStamp, TransactionStatus," paragraph.
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$xml10 = '<w:p ' + $wns + ' w:rsidR="00CB0198" w:rsidRDefault="00CB0198"><w:r><w:tab/><w:t xml:space="preserve">TransactionId, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TimeStamp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, TransactionStatus,</w:t></w:r></w:p>'
$p10.Range.InsertXML($xml10) | Out-Null

# ---------------------------------------------------------------------------
# 3) Wrap "RateStars" with spellcheck proofErr markers.
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$xml12 = '<w:p ' + $wns + ' w:rsidR="00376909" w:rsidRDefault="00CB0198" w:rsidP="00376909"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>RateStars</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Comments</w:t></w:r></w:p>'
$p12.Range.InsertXML($xml12) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the redundant "CourseId -> EmailAddress, ..." FD line, leaving
#    just a lone tab character in its paragraph.
# ---------------------------------------------------------------------------
$p23 = $d.Paragraphs(23)
$xml23 = '<w:p ' + $wns + ' w:rsidR="00BF5764" w:rsidRDefault="00BF5764" w:rsidP="00E40EB1"><w:r><w:tab/></w:r></w:p>'
$p23.Range.InsertXML($xml23) | Out-Null

# ---------------------------------------------------------------------------
# 5) Rewrite the R1-R7 relation schema paragraphs, collapsing the
#    multi-run text into the new run layout and dropping the bookmark
#    (it was relocated to the EmailAddress split above).
# ---------------------------------------------------------------------------
$p27 = $d.Paragraphs(27)
$xml27 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:t>R1(Email</w:t></w:r><w:r><w:t>Address, StudentName, BirthDate, Location, StudentGender)</w:t></w:r></w:p>'
$p27.Range.InsertXML($xml27) | Out-Null

$p28 = $d.Paragraphs(28)
$xml28 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:t xml:space="preserve">R2(CourseId, </w:t></w:r><w:r><w:t>CourseName, Overview, Duration, Difficulty, Category, FAQ)</w:t></w:r></w:p>'
$p28.Range.InsertXML($xml28) | Out-Null

$p29 = $d.Paragraphs(29)
$xml29 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRPr="00AD6414" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:t>R3(UniversityId, UniversityName, Description)</w:t></w:r></w:p>'
$p29.Range.InsertXML($xml29) | Out-Null

$p30 = $d.Paragraphs(30)
$xml30 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:lastRenderedPageBreak/><w:t>R4(InstructorId -&gt; InstructorName, InstructorGender, Specialization)</w:t></w:r></w:p>'
$p30.Range.InsertXML($xml30) | Out-Null

$p31 = $d.Paragraphs(31)
$xml31 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:t>R5(TransactionId, Timestamp, TransactionStatus)</w:t></w:r></w:p>'
$p31.Range.InsertXML($xml31) | Out-Null

$p32 = $d.Paragraphs(32)
$xml32 = '<w:p ' + $wns + ' w:rsidR="0091451D" w:rsidRDefault="0091451D" w:rsidP="0091451D"><w:r><w:t>R6(CreditCardNo, NameOnCard, CVV, ExpiryDate)</w:t></w:r></w:p>'
$p32.Range.InsertXML($xml32) | Out-Null

$p33 = $d.Paragraphs(33)
$xml33 = '<w:p ' + $wns + ' w:rsidR="00BF5764" w:rsidRDefault="00BF5764" w:rsidP="00BF5764"><w:r><w:t>R7(CourseId, EmailAddress, UniversityId, InstructorId, TransactionId, CreditCardNo)</w:t></w:r></w:p>'
$p33.Range.InsertXML($xml33) | Out-Null
